# UI bug fixed (navigation removed), healthbars added!
#
# Adds three new paragraphs describing the XP / leveling / leaderboard
# design right after the "increases amount of load you can take with
# you" bullet, and moves the (hidden) "_GoBack" bookmark from that
# bullet to the end of the final new paragraph.

$d = $word.ActiveDocument

# The bullet paragraph the new content goes after.
$anchor = $d.Paragraphs(9)

# Make room: three fresh paragraphs right after the anchor bullet.
$anchor.Range.InsertParagraphAfter()
$anchor.Range.InsertParagraphAfter()
$anchor.Range.InsertParagraphAfter()

# --- Paragraph 1: "The player is rewarded xp for his actions..." ---
$r = $d.Paragraphs(10).Range
$r.Collapse(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The player is rewarded </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>xp</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> for his actions, killing enemies, acquiring resources etc.</w:t></w:r>' +
  '</w:p>'
$r.InsertXML($xml)

# --- Paragraph 2: "The goal is to achieve a certain amount of average xp..." ---
$r = $d.Paragraphs(11).Range
$r.Collapse(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The goal is to achieve a certain amount of average </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>xp</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> over a few consecutive runs. This will unlock a teleport upgrade that will unlock new area types, enemies and other content.</w:t></w:r>' +
  '</w:p>'
$r.InsertXML($xml)

# --- Paragraph 3: "After reaching the maximum level..." plus the _GoBack bookmark ---
$r = $d.Paragraphs(12).Range
$r.Collapse(1)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>After reaching the maximum level with a character, the goal becomes to compete against other players by being ranked in leaderboards.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$r.InsertXML($xml)

# The _GoBack bookmark used to sit at the end of the anchor bullet;
# rewrite that paragraph without it now that it has moved above.
$r = $anchor.Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="008F28C7" w:rsidRDefault="008F28C7" w:rsidP="008F28C7">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>increases amount of load you can take with you</w:t></w:r>' +
  '</w:p>'
$r.InsertXML($xml)
